# Delete the entire spreadsheet row 75 ("ドッキリカメラ" post), which shifts
# all subsequent rows up by one. This matches the commit's net effect:
# dimension shrinks from A1:C238 to A1:C237.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(75).Delete()
